# Updated symbol list on Fri Dec 16 15:26:28 UTC 2022 with GitHub Actions
#
# Refresh the live crypto price snapshot on Sheet1. Numeric-looking values
# are stored as text in this sheet (matching the original inline-string
# layout), so each numeric update is written with a leading apostrophe to
# keep Excel from re-typing the cell as a Number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) refresh for existing rows -----------------------
$ws.Range("D2").Value  = "'248.69"
$ws.Range("D3").Value  = "'24.17"
$ws.Range("D4").Value  = "'5.961"
$ws.Range("D5").Value  = "'0.05888"
$ws.Range("D6").Value  = "'3.432"
$ws.Range("D7").Value  = "'6.523"
$ws.Range("D8").Value  = "'1.328"
$ws.Range("D9").Value  = "'0.7967"
$ws.Range("D10").Value = "'0.1474"
$ws.Range("D11").Value = "'0.07753"
$ws.Range("D12").Value = "'0.03299"
$ws.Range("D13").Value = "'0.03015"
$ws.Range("D14").Value = "'0.09213"
$ws.Range("D15").Value = "'3.577"
$ws.Range("D16").Value = "'0.001684"
$ws.Range("D17").Value = "'0.04782"

$ws.Range("D18").Value = "'0.0006038"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("D19").Value = "'0.006224"
$ws.Range("D20").Value = "'0.005550"
$ws.Range("D21").Value = "'0.001070"
$ws.Range("D22").Value = "'0.0001502"
$ws.Range("D23").Value = "'3.697"
$ws.Range("D24").Value = "'2.209"
$ws.Range("D25").Value = "'0.3351"
$ws.Range("D26").Value = "'0.1254"
$ws.Range("D27").Value = "'0.0006277"

$ws.Range("D40").Value = "'0.04387"
$ws.Range("D41").Value = "'0.007028"

# --- Rows 42/43 swapped positions in the ranking -----------------------
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1065"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003225"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "'0.009645"
$ws.Range("D45").Value = "'0.002463"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D48").Value = "'0.9912"

$ws.Range("D49").Value = "'0.1110"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"

$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D51").Value = "'0.01011"
